# The deck's "Integral" design (ppt/theme/theme1.xml, used by the slide
# master / all slides) is being swapped for the stock "Office Theme"
# palette that previously lived only in the notes-master theme
# (ppt/theme/theme2.xml). Re-point every themed colour slot on the
# slide master's theme to the Office Theme RGB values.
$p = $ppt.ActivePresentation
$s1 = $p.Slides.Item(1)
$tcs = $s1.ThemeColorScheme

# RGB() packs as 0x00BBGGRR, so build each value from the target hex.
function RgbVal([int]$r, [int]$g, [int]$b) {
    return $r + ($g * 256) + ($b * 65536)
}

$tcs.Colors(3).RGB  = (RgbVal 0x44 0x54 0x6A)   # dk2       -> 44546A
$tcs.Colors(4).RGB  = (RgbVal 0xE7 0xE6 0xE6)   # lt2       -> E7E6E6
$tcs.Colors(5).RGB  = (RgbVal 0x5B 0x9B 0xD5)   # accent1   -> 5B9BD5
$tcs.Colors(6).RGB  = (RgbVal 0xED 0x7D 0x31)   # accent2   -> ED7D31
$tcs.Colors(7).RGB  = (RgbVal 0xA5 0xA5 0xA5)   # accent3   -> A5A5A5
$tcs.Colors(8).RGB  = (RgbVal 0xFF 0xC0 0x00)   # accent4   -> FFC000
$tcs.Colors(9).RGB  = (RgbVal 0x44 0x72 0xC4)   # accent5   -> 4472C4
$tcs.Colors(10).RGB = (RgbVal 0x70 0xAD 0x47)   # accent6   -> 70AD47
$tcs.Colors(11).RGB = (RgbVal 0x05 0x63 0xC1)   # hlink     -> 0563C1
$tcs.Colors(12).RGB = (RgbVal 0x95 0x4F 0x72)   # folHlink  -> 954F72
